# Update "vendas_atipicas" sheet with refreshed atypical-sales data.
# The table is fully reloaded with a new date window (2025-07-16 .. 2025-07-30),
# adding one extra row (row 11) and refreshing the numeric figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 11

# Columns A, C, D, F hold text values in this sheet (dates and id-like
# strings are kept as text, not auto-converted numbers/dates), so force
# a text number format on them before writing, then drop back to the
# default "Normal" style once done so no visible style attribute remains
# on the data cells.
$textCols = @("A", "C", "D", "F")
foreach ($col in $textCols) {
    $ws.Range("$col$firstRow`:$col$lastRow").NumberFormat = "@"
}

# New table contents, row => (A data, B quantidade_atipica, C cliente,
# D id_venda, E id_produto, F produto, G estoque_atualizado,
# H media_vendas, I desvio_padrao)
$rows = @{
    2  = @("2025-07-16", 4, "BEMOL S/A", "383665", 3984,  "BARALHO PLASTICO 1001 COPAG ESTOJO C/2 110 UNIDADES",          -37,  1.19, 0.68)
    3  = @("2025-07-17", 2, "BEMOL S/A", "384275", 13185, "KIT LANCHE FUNDO DO MAR GARRAFA PLASTICA 500ML + MARMITA 700ML", 2,  1.03, 0.17)
    4  = @("2025-07-21", 2, "BEMOL S/A", "386126", 10114, "CARREGADOR USB-C A GOLD 20W CA31-4",                            -94,  1.05, 0.22)
    5  = @("2025-07-21", 2, "BEMOL S/A", "386260", 47869, "MARMITA ELÉTRICA ONEX",                                        -17,  1.08, 0.29)
    6  = @("2025-07-28", 2, "BEMOL S/A", "387489", 13242, "POWER BANK 5000mAh PN-952",                                   -245,  1.03, 0.16)
    7  = @("2025-07-29", 2, "BEMOL S/A", "389107", 10525, "BOMBA AUTOMATICA PARA GALAO DE AGUA RECARREGAVEL USB",        -103,  1.04, 0.19)
    8  = @("2025-07-29", 2, "BEMOL S/A", "389675", 49177, "KIT SMARTWATCH INOVA PULSEIRA PRETO LISA (FONE+FONTE+CABO)",   -65,  1.03, 0.18)
    9  = @("2025-07-29", 2, "BEMOL S/A", "389699", 48696, "FONE DE OUVIDO TIPO-C EJ-105",                                 -73,  1.06, 0.23)
    10 = @("2025-07-30", 2, "BEMOL S/A", "390273", 13958, "CAMPAINHA RESIDENCIAL SEM FIO 32 TOQUES C/PLUG BIVOL",         -25,  1.08, 0.27)
    11 = @("2025-07-30", 2, "BEMOL S/A", "390286", 12945, "FONE DE OUVIDO SEM FIO BT BASIKE FON-9856",                    -50,  1.03, 0.18)
}

foreach ($r in $firstRow..$lastRow) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 1).Value = $vals[0]   # A data
    $ws.Cells.Item($r, 2).Value = $vals[1]   # B quantidade_atipica
    $ws.Cells.Item($r, 3).Value = $vals[2]   # C cliente
    $ws.Cells.Item($r, 4).Value = $vals[3]   # D id_venda
    $ws.Cells.Item($r, 5).Value = $vals[4]   # E id_produto
    $ws.Cells.Item($r, 6).Value = $vals[5]   # F produto
    $ws.Cells.Item($r, 7).Value = $vals[6]   # G estoque_atualizado
    $ws.Cells.Item($r, 8).Value = $vals[7]   # H media_vendas
    $ws.Cells.Item($r, 9).Value = $vals[8]   # I desvio_padrao
}

# Drop the temporary text formatting back to the default style so the
# data rows end up without any explicit style index, matching the rest
# of the sheet.
foreach ($col in $textCols) {
    $ws.Range("$col$firstRow`:$col$lastRow").Style = "Normal"
}

Write-Output "vendas_atipicas updated: rows 2-11 refreshed"
